# Fruta / hortaliza, semanal
# Insert a new weekly record as row 14 in the "Tuna" sheet, pushing the
# existing rows 14-99 down to 15-100 (dimension grows from A1:T99 to A1:T100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 14; Excel copies row 13's
# formatting (including the date number format on column D) down onto it.
$ws.Rows.Item(14).Insert()

# Populate the new row with the new daily observation. Columns that are
# identical across every record in this sheet (market/region/product
# metadata) are copied straight from the neighbouring rows.
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 45061
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100107
$ws.Cells.Item(14, 8).Value = "Otros"
$ws.Cells.Item(14, 9).Value = 100107011
$ws.Cells.Item(14, 10).Value = "Tuna"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 20000
$ws.Cells.Item(14, 15).Value = 20000
$ws.Cells.Item(14, 16).Value = 20000
$ws.Cells.Item(14, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(14, 19).Value = 1250
$ws.Cells.Item(14, 20).Value = 16
